# "Generate Report for Handback" - fills in the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns on the per-language sheets
# once a handback has been processed, and flips the Overview sheet's status
# text to reflect that the content is now in sync with en-US.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet: status cells move from "In Translation" to the
#    "handed back" message.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# widen the zh-cn / de-de status columns so the longer text is readable
$overview.Columns.Item(5).ColumnWidth = 29.16666667
$overview.Columns.Item(6).ColumnWidth = 29.16666667

# ---------------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = 29.16666667
$zhcn.Columns.Item(9).ColumnWidth = 39.16666667
$zhcn.Columns.Item(10).ColumnWidth = 39.16666667

# Row 2 - 16a7712c-38af-47cf-b77b-bc11c6d12cc8
$zhcn.Range("I2").Value = "16a7712c-38af-47cf-b77b-bc11c6d12cc8.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9eab9ea91f62f6b13f02656e7503917cf41ae92f/e2e/16a7712c-38af-47cf-b77b-bc11c6d12cc8.md", "", "", "16a7712c-38af-47cf-b77b-bc11c6d12cc8.md")
$zhcn.Range("J2").Value = "16a7712c-38af-47cf-b77b-bc11c6d12cc8.7b49f57cbf8e2ea8f2fc987491980ddbe88b6823.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-12 12:26:12"

# Row 3 - 7424f250-d920-4364-a36b-37196aaf66a2
$zhcn.Range("I3").Value = "7424f250-d920-4364-a36b-37196aaf66a2.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9eab9ea91f62f6b13f02656e7503917cf41ae92f/e2e/7424f250-d920-4364-a36b-37196aaf66a2.md", "", "", "7424f250-d920-4364-a36b-37196aaf66a2.md")
$zhcn.Range("J3").Value = "7424f250-d920-4364-a36b-37196aaf66a2.c4c38fd7327a428b239766686ad15505a814c5d5.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-12 12:26:12"

$zhcn.Range("I2").Style = "Hyperlink"
$zhcn.Range("I3").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = 29.16666667
$dede.Columns.Item(9).ColumnWidth = 39.16666667
$dede.Columns.Item(10).ColumnWidth = 39.16666667

# Row 2 - 16a7712c-38af-47cf-b77b-bc11c6d12cc8
$dede.Range("I2").Value = "16a7712c-38af-47cf-b77b-bc11c6d12cc8.md"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9eab9ea91f62f6b13f02656e7503917cf41ae92f/e2e/16a7712c-38af-47cf-b77b-bc11c6d12cc8.md", "", "", "16a7712c-38af-47cf-b77b-bc11c6d12cc8.md")
$dede.Range("J2").Value = "16a7712c-38af-47cf-b77b-bc11c6d12cc8.7b49f57cbf8e2ea8f2fc987491980ddbe88b6823.de-de.xlf"
$dede.Range("K2").Value = "2016-08-12 12:26:21"

# Row 3 - 7424f250-d920-4364-a36b-37196aaf66a2
$dede.Range("I3").Value = "7424f250-d920-4364-a36b-37196aaf66a2.md"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9eab9ea91f62f6b13f02656e7503917cf41ae92f/e2e/7424f250-d920-4364-a36b-37196aaf66a2.md", "", "", "7424f250-d920-4364-a36b-37196aaf66a2.md")
$dede.Range("J3").Value = "7424f250-d920-4364-a36b-37196aaf66a2.c4c38fd7327a428b239766686ad15505a814c5d5.de-de.xlf"
$dede.Range("K3").Value = "2016-08-12 12:26:21"

$dede.Range("I2").Style = "Hyperlink"
$dede.Range("I3").Style = "Hyperlink"
